$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TestBean -> JavaBean throughout the two affected cells.
# F4 ("Data TestBean beans2") is edited first, then B3
# ("Method String print(TestBean bean)") -- matching the original authoring
# order so the shared-string table comes out in the same slot order.
$ws.Range("F4").Value2 = "Data JavaBean beans2"
$ws.Range("B3").Value2 = "Method String print(JavaBean bean)"

# Update the selection left behind in the saved sheet view.
$ws.Range("E26").Select()

# Rename the default cell style from the Russian-locale "Обычный" to the
# canonical "Normal". The workbook carries a single (builtin) cell style;
# deleting it collapses it back to the canonical default name.
$styles = $wb.Styles
$styles.Item(1).Delete()

# Localize the theme part's display names (Russian defaults -> English
# defaults), using the documented Excel Theme object model.
$theme = $wb.Theme
$theme.Name = "Office Theme"
$theme.ThemeColorScheme.Name = "Office"
$theme.ThemeFontScheme.Name = "Office"
